# Apply updated transition-matrix probabilities to Sheet1.
# These values represent the simulated game-state transition matrix
# after adding more games / speeding up the simulation logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2"  = 0.1428571428571428
    "C2"  = 0.5714285714285714
    "P2"  = 0.1428571428571428
    "S2"  = 0.1428571428571428

    "P3"  = 1

    "B6"  = 0.2222222222222222
    "J6"  = 0.2222222222222222
    "Q6"  = 0.1111111111111111
    "S6"  = 0.4444444444444444

    "Q7"  = 0.2857142857142857
    "S7"  = 0.7142857142857143

    "B8"  = 0.06451612903225806
    "F8"  = 0.03225806451612903
    "J8"  = 0.1612903225806452
    "Q8"  = 0.09677419354838709
    "S8"  = 0.6451612903225806

    "B9"  = 0.09090909090909091
    "Q9"  = 0.1818181818181818
    "S9"  = 0.7272727272727273

    "B10" = 0.02380952380952381
    "F10" = 0.119047619047619
    "J10" = 0.119047619047619
    "O10" = 0.02380952380952381
    "Q10" = 0.1666666666666667
    "R10" = 0.04761904761904762
    "S10" = 0.5

    "G11" = 0.2727272727272727
    "K11" = 0.09090909090909091
    "L11" = 0.6363636363636364

    "G12" = 0.5714285714285714
    "J12" = 0.1428571428571428
    "K12" = 0.2857142857142857

    "J13" = 1

    "F15" = 0.2
    "H15" = 0.2
    "J15" = 0.2
    "K15" = 0.2
    "M15" = 0.2

    "I16" = 0.2
    "J16" = 0.2
    "K16" = 0.2
    "S16" = 0.4

    "H17" = 0.3333333333333333
    "I17" = 0.06666666666666667
    "J17" = 0.4666666666666667
    "K17" = 0.06666666666666667
    "S17" = 0.06666666666666667

    "O18" = 0.5
    "S18" = 0.5

    "F19" = 0.01388888888888889
    "H19" = 0.3333333333333333
    "I19" = 0.1388888888888889
    "J19" = 0.2638888888888889
    "K19" = 0.06944444444444445
    "N19" = 0.01388888888888889
    "O19" = 0.02777777777777778
    "S19" = 0.1388888888888889
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
